# Previo a entrega minuta 29-ago
# Update status cells: two previously-"Pendiente ..." statuses are now resolved ("OK"),
# and the now-unused shared strings are thereby dropped from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Actividad 3): "Pendiente agregar tipo_asociado en la respuesta del ws" -> "OK"
$ws.Range("C7").Value = "OK"

# Row 14 (Actividad 10): "Pendiente codigos validos para mascotas " -> "OK"
$ws.Range("C14").Value = "OK"

# Update the active selection to reflect where the editor left off
$ws.Range("C22").Select()
